$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes the existing rows 3..5 down to 4..6)
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly price record
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44664
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107011
$ws.Range("J3").Value = "Tuna"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 1639
$ws.Range("T3").Value = 18
